# Vector_Comp.py - Fixed error in calculation of scale factor
# (2d & 3d vector magnitude calculation indices were wrong)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Points_2d")

# Correct the scale-factor inputs that were left blank due to the
# wrong indices used in the magnitude calculation.
$ws.Range("C8").Value = 3.5
$ws.Range("D8").Value = 3.5
$ws.Range("A9").Value = 1.25
$ws.Range("B9").Value = 1.25

# Update the active selection / view to match the author's final state.
$ws.Activate()
$ws.Range("E8").Select()
